# Custo de reposição.xlsx - apply PIS formula fix (divide by 100) and
# adjust the number formats for the PIS (G) and COFINS (H) columns of the
# "DADOS" table so that the new PIS values (now a fraction instead of a
# percent-like number) display with more decimal places.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix the PIS formula (column G) for every data row of the table: the
#    previous formula produced a number 100x too large, so we divide by 100.
for ($r = 2; $r -le 127; $r++) {
    $ws.Cells.Item($r, 7).Formula = "=1.65 * (1-(DADOS[[#This Row],[ICMS]]))/100"
}

# 2) Give the PIS column (G) a number format with more decimals so the
#    now-smaller values remain legible.
$ws.Range("G2:G127").NumberFormat = "0.00000"

# 3) Give the COFINS column (H) - including its header - a matching number
#    format with extra decimals as well.
$ws.Range("H1:H127").NumberFormat = "0.0000"
